$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename: "Done?" -> "Read?" ---
$ws.Range("D1").Value = "Read?"

# --- Mark a couple of existing rows with an X in the (now "Read?") column ---
$ws.Range("D3").Value = "X"
$ws.Range("D11").Value = "X"

# --- Rename " - Procedure" to " - Procedure(s)" ---
$ws.Range("A22").Value = " - Procedure(s)"

# --- Give the (still empty) Results date cell the same date formatting as its neighbours ---
$ws.Range("C23").NumberFormat = "d-mmm"

# --- New "Extra Work" section ---
$ws.Range("A26").Value = "Extra Work"
$ws.Range("A26").Font.Italic = $true

$ws.Range("A27").Value = "Introduction fix table"
$ws.Range("B27").Value = "Nicklas"
$ws.Range("C27").Value = "2017-11-03"
$ws.Range("C27").NumberFormat = "d-mmm"

$ws.Range("A28").Value = "Source on price of RAMIS"
$ws.Range("B28").Value = "Oliver"
$ws.Range("C28").Value = "2017-11-03"
$ws.Range("C28").NumberFormat = "d-mmm"

$ws.Range("A29").Value = "Include tools in Artefact model"
$ws.Range("B29").Value = "Stjernholm"
$ws.Range("C29").Value = "2017-11-03"
$ws.Range("C29").NumberFormat = "d-mmm"

$ws.Range("A30").Value = "Merge interviews and obs."
$ws.Range("B30").Value = "Freddie"
$ws.Range("C30").Value = "2017-11-03"
$ws.Range("C30").NumberFormat = "d-mmm"

$ws.Range("A31").Value = "Interview conclusions"
$ws.Range("B31").Value = "Freddie"
$ws.Range("C31").Value = "2017-11-03"
$ws.Range("C31").NumberFormat = "d-mmm"

$ws.Range("A32").Value = "Add Johan's interview"
$ws.Range("B32").Value = "Freddie"
$ws.Range("C32").Value = "2017-11-03"
$ws.Range("C32").NumberFormat = "d-mmm"

$ws.Range("A33").Value = "Rewrite rules"
$ws.Range("B33").Value = "Atanas"
$ws.Range("C33").Value = "2017-11-03"
$ws.Range("C33").NumberFormat = "d-mmm"

# --- Column A best-fit width (new, longer task names no longer fit) ---
$ws.Columns.Item(1).AutoFit()

# --- Selection as left by the editing author ---
$ws.Range("A6").Select()
